$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027016220519477
$ws.Range("D2").Value = 1.057198492128738
$ws.Range("E2").Value = 1.02716981780513
$ws.Range("F2").Value = 1.06069212986105
$ws.Range("I2").Value = 1.045288073939175
$ws.Range("J2").Value = 1.032176441771654
$ws.Range("K2").Value = 1.059934118266941
$ws.Range("L2").Value = 1.029990328853348
$ws.Range("M2").Value = 1.063418218300183
$ws.Range("N2").Value = 1.014752633434396

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02781036629354
$ws.Range("D3").Value = 1.057863918442534
$ws.Range("E3").Value = 1.02784004761181
$ws.Range("F3").Value = 1.061480388565241
$ws.Range("I3").Value = 1.045510797108564
$ws.Range("J3").Value = 1.032611392166664
$ws.Range("K3").Value = 1.060413710810984
$ws.Range("L3").Value = 1.03046859913562
$ws.Range("M3").Value = 1.064021018302289
$ws.Range("N3").Value = 1.014897499079952

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028324893595774
$ws.Range("D4").Value = 1.05829475110635
$ws.Range("E4").Value = 1.028274694079999
$ws.Range("F4").Value = 1.06199112221181
$ws.Range("I4").Value = 1.045653839276417
$ws.Range("J4").Value = 1.032892827445421
$ws.Range("K4").Value = 1.060723612029304
$ws.Range("L4").Value = 1.030778353563444
$ws.Range("M4").Value = 1.064411087465663
$ws.Range("N4").Value = 1.014991206390631

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028541357584426
$ws.Range("D5").Value = 1.05847593274827
$ws.Range("E5").Value = 1.028457648120471
$ws.Range("F5").Value = 1.06220599471395
$ws.Range("I5").Value = 1.045713716003206
$ws.Range("J5").Value = 1.033011139973656
$ws.Range("K5").Value = 1.060853790637942
$ws.Range("L5").Value = 1.030908640350359
$ws.Range("M5").Value = 1.064575074884002
$ws.Range("N5").Value = 1.015030593200117

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028577711979497
$ws.Range("D6").Value = 1.058506357378826
$ws.Range("E6").Value = 1.028488380263488
$ws.Range("F6").Value = 1.062242082077538
$ws.Range("I6").Value = 1.045723754400913
$ws.Range("J6").Value = 1.033031004963322
$ws.Range("K6").Value = 1.060875642062436
$ws.Range("L6").Value = 1.030930519906667
$ws.Range("M6").Value = 1.064602609178423
$ws.Range("N6").Value = 1.015037205951411

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028327785384602
$ws.Range("D7").Value = 1.05829717183271
$ws.Range("E7").Value = 1.0282771378247
$ws.Range("F7").Value = 1.061993992722387
$ws.Range("I7").Value = 1.045654640367893
$ws.Range("J7").Value = 1.03289440835604
$ws.Range("K7").Value = 1.060725351892349
$ws.Range("L7").Value = 1.030780094204172
$ws.Range("M7").Value = 1.064413278666619
$ws.Range("N7").Value = 1.014991732710107

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027284467704724
$ws.Range("D8").Value = 1.057423321569822
$ws.Range("E8").Value = 1.027396124593797
$ws.Range("F8").Value = 1.060958384007747
$ws.Range("I8").Value = 1.04536356622085
$ws.Range("J8").Value = 1.032323435928546
$ws.Range("K8").Value = 1.060096286320646
$ws.Range("L8").Value = 1.030151903401049
$ws.Range("M8").Value = 1.063621932735808
$ws.Range("N8").Value = 1.014801597489162

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025451155141479
$ws.Range("D9").Value = 1.055885545858284
$ws.Range("E9").Value = 1.025851124558645
$ws.Range("F9").Value = 1.059138792584699
$ws.Range("I9").Value = 1.044842467482162
$ws.Range("J9").Value = 1.031317313695253
$ws.Range("K9").Value = 1.058984590373812
$ws.Range("L9").Value = 1.029047168630829
$ws.Range("M9").Value = 1.062227691282418
$ws.Range("N9").Value = 1.014466341813221

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024232511116412
$ws.Range("D10").Value = 1.054861867514172
$ws.Range("E10").Value = 1.024826250496679
$ws.Range("F10").Value = 1.057929411059998
$ws.Range("I10").Value = 1.044489618352516
$ws.Range("J10").Value = 1.030646643253255
$ws.Range("K10").Value = 1.058241401340725
$ws.Range("L10").Value = 1.028312252426434
$ws.Range("M10").Value = 1.061298445629981
$ws.Range("N10").Value = 1.014242721836883

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023705692176351
$ws.Range("D11").Value = 1.054418986012215
$ws.Range("E11").Value = 1.024383708190927
$ws.Range("F11").Value = 1.057406635091764
$ws.Range("I11").Value = 1.044335549876282
$ws.Range("J11").Value = 1.030356269803785
$ws.Range("K11").Value = 1.057919124003223
$ws.Range("L11").Value = 1.027994416310727
$ws.Range("M11").Value = 1.060896153509654
$ws.Range("N11").Value = 1.014145870084482

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023510139436089
$ws.Range("D12").Value = 1.054254538949533
$ws.Range("E12").Value = 1.024219515737356
$ws.Range("F12").Value = 1.057212589343784
$ws.Range("I12").Value = 1.04427813023854
$ws.Range("J12").Value = 1.030248418150407
$ws.Range("K12").Value = 1.057799346604308
$ws.Range("L12").Value = 1.027876417450795
$ws.Range("M12").Value = 1.060746737747348
$ws.Range("N12").Value = 1.014109892063665

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023552080193977
$ws.Range("D13").Value = 1.054289810706745
$ws.Range("E13").Value = 1.024254727061977
$ws.Range("F13").Value = 1.057254206602525
$ws.Range("I13").Value = 1.044290455607687
$ws.Range("J13").Value = 1.030271552405332
$ws.Range("K13").Value = 1.057825042370785
$ws.Range("L13").Value = 1.027901725885944
$ws.Range("M13").Value = 1.060778787298063
$ws.Range("N13").Value = 1.014117609598143

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023689525037711
$ws.Range("D14").Value = 1.054405391554365
$ws.Range("E14").Value = 1.024370132156163
$ws.Range("F14").Value = 1.057390592412085
$ws.Range("I14").Value = 1.044330807457279
$ws.Range("J14").Value = 1.030347354623218
$ws.Range("K14").Value = 1.057909224569402
$ws.Range("L14").Value = 1.027984661263645
$ws.Range("M14").Value = 1.060883802471387
$ws.Range("N14").Value = 1.014142896187928

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023774226789814
$ws.Range("D15").Value = 1.054476612588509
$ws.Range("E15").Value = 1.024441261941848
$ws.Range("F15").Value = 1.057474642366407
$ws.Range("I15").Value = 1.044355644183295
$ws.Range("J15").Value = 1.030394059701282
$ws.Range("K15").Value = 1.057961082872451
$ws.Range("L15").Value = 1.028035768430257
$ws.Range("M15").Value = 1.060948507621314
$ws.Range("N15").Value = 1.014158475711019

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024267492666093
$ws.Range("D16").Value = 1.054891268240288
$ws.Range("E16").Value = 1.024855646755275
$ws.Range("F16").Value = 1.057964125005281
$ws.Range("I16").Value = 1.044499816429102
$ws.Range("J16").Value = 1.030665915150655
$ws.Range("K16").Value = 1.05826278002744
$ws.Range("L16").Value = 1.02833335444663
$ws.Range("M16").Value = 1.061325146218597
$ws.Range("N16").Value = 1.014249149127214

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024577137339833
$ws.Range("D17").Value = 1.055151473425713
$ws.Range("E17").Value = 1.025115911388693
$ws.Range("F17").Value = 1.058271405477013
$ws.Range("I17").Value = 1.044589909172407
$ws.Range("J17").Value = 1.030836452166316
$ws.Range("K17").Value = 1.058451901558609
$ws.Range("L17").Value = 1.028520127023557
$ws.Range("M17").Value = 1.061561423660502
$ws.Range("N17").Value = 1.014306020380414

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024757830887521
$ws.Range("D18").Value = 1.055303283039332
$ws.Range("E18").Value = 1.02526783836932
$ws.Range("F18").Value = 1.058450723110034
$ws.Range("I18").Value = 1.044642334942608
$ws.Range("J18").Value = 1.030935926454642
$ws.Range("K18").Value = 1.058562167339796
$ws.Range("L18").Value = 1.028629105596293
$ws.Range("M18").Value = 1.06169924763903
$ws.Range("N18").Value = 1.014339190186808

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024819456750341
$ws.Range("D19").Value = 1.055355052294165
$ws.Range("E19").Value = 1.025319661673236
$ws.Range("F19").Value = 1.058511880332207
$ws.Range("I19").Value = 1.044660189748024
$ws.Range("J19").Value = 1.030969845074873
$ws.Range("K19").Value = 1.0585997573369
$ws.Range("L19").Value = 1.028666270755048
$ws.Range("M19").Value = 1.061746243244383
$ws.Range("N19").Value = 1.014350499840718

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024543906793947
$ws.Range("D20").Value = 1.055123552099097
$ws.Range("E20").Value = 1.025087975131594
$ws.Range("F20").Value = 1.058238428265653
$ws.Range("I20").Value = 1.044580255872455
$ws.Range("J20").Value = 1.030818154845974
$ws.Range("K20").Value = 1.058431615314049
$ws.Range("L20").Value = 1.028500084226512
$ws.Range("M20").Value = 1.061536072554335
$ws.Range("N20").Value = 1.01429991886085

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023649047319015
$ws.Range("D21").Value = 1.054371354234543
$ws.Range("E21").Value = 1.02433614304195
$ws.Range("F21").Value = 1.057350426404309
$ws.Range("I21").Value = 1.044318930128955
$ws.Range("J21").Value = 1.030325032584203
$ws.Range("K21").Value = 1.057884436908764
$ws.Range("L21").Value = 1.027960237204723
$ws.Range("M21").Value = 1.060852877734588
$ws.Range("N21").Value = 1.014135449995704

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023087174019719
$ws.Range("D22").Value = 1.053898759066943
$ws.Range("E22").Value = 1.023864521922457
$ws.Range("F22").Value = 1.056792895346937
$ws.Range("I22").Value = 1.044153515341831
$ws.Range("J22").Value = 1.030015022244459
$ws.Range("K22").Value = 1.05754000428239
$ws.Range("L22").Value = 1.027621159627212
$ws.Range("M22").Value = 1.060423404712491
$ws.Range("N22").Value = 1.014032024950715

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023384961020789
$ws.Range("D23").Value = 1.054149257651408
$ws.Range("E23").Value = 1.024114433659879
$ws.Range("F23").Value = 1.057088377409067
$ws.Range("I23").Value = 1.044241309651796
$ws.Range("J23").Value = 1.030179360829056
$ws.Range("K23").Value = 1.057722631900231
$ws.Range("L23").Value = 1.027800877828907
$ws.Range("M23").Value = 1.060651068423185
$ws.Range("N23").Value = 1.014086853975745

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024558921990859
$ws.Range("D24").Value = 1.055136168432744
$ws.Range("E24").Value = 1.025100597956171
$ws.Range("F24").Value = 1.058253328982124
$ws.Range("I24").Value = 1.044584618165825
$ws.Range("J24").Value = 1.030826422608541
$ws.Range("K24").Value = 1.058440781935879
$ws.Range("L24").Value = 1.028509140588734
$ws.Range("M24").Value = 1.061547527605526
$ws.Range("N24").Value = 1.014302675881865

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025924489774136
$ws.Range("D25").Value = 1.056282841192611
$ws.Range("E25").Value = 1.026249648784016
$ws.Range("F25").Value = 1.059608561358852
$ws.Range("I25").Value = 1.044978148678301
$ws.Range("J25").Value = 1.031577412195258
$ws.Range("K25").Value = 1.059272360141494
$ws.Range("L25").Value = 1.029332497725377
$ws.Range("M25").Value = 1.062588100094642
$ws.Range("N25").Value = 1.014553035903381
